$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add three new rows of data (16, 17, 18) continuing the date sequence
$ws.Range("A16").Value = 44230
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0

$ws.Range("A17").Value = 44231
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0

$ws.Range("A18").Value = 44232
$ws.Range("B18").Value = 24
$ws.Range("C18").Value = 0

# Match the date number format used by the rest of column A by copying
# the formatting (not just the raw number format string) from an existing cell
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update selection to match the diff (active cell C18)
$ws.Range("C18").Select()
